$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement dataset (rows 2-11), columns: Employee ID, Name, Meal Type, Time, Date, City
$data = @(
    @(56123, "Ava Martinez",   "Normal", "2:30 PM", "18 July 2024", "Lahore"),
    @(18392, "Hina Saeed",     "Normal", "3:00 PM", "17 July 2024", "Lahore"),
    @(16254, "David Williams", "Normal", "1:30 PM", "18 July 2024", "Lahore"),
    @(21789, "Sophia Wilson",  "Normal", "2:00 PM", "17 July 2024", "Lahore"),
    @(23781, "Fatima Yousaf",  "Normal", "3:00 PM", "19 July 2024", "Lahore"),
    @(74528, "Sana Abbas",     "Normal", "1:30 PM", "15 July 2024", "Lahore"),
    @(74528, "Sana Abbas",     "Normal", "1:30 PM", "17 July 2024", "Lahore"),
    @(74528, "Sana Abbas",     "Normal", "1:30 PM", "18 July 2024", "Lahore"),
    @(90432, "Ayesha Ahmed",   "Diet",   "2:30 PM", "18 July 2024", "Lahore"),
    @(90432, "Ayesha Ahmed",   "Normal", "1:30 PM", "18 July 2024", "Lahore")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
